$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix existing row 15, column A value (tiny float precision correction)
$ws.Cells.Item(15, 1).Value = 45864.79193502315

# Add new row 16
$ws.Cells.Item(16, 1).Value = 45864.83358132153
$ws.Cells.Item(16, 1).NumberFormat = $ws.Cells.Item(15, 1).NumberFormat
$ws.Cells.Item(16, 2).Value = 2025
$ws.Cells.Item(16, 3).Value = 30
$ws.Cells.Item(16, 4).Value = 13.56
$ws.Cells.Item(16, 5).Value = 88.73
$ws.Cells.Item(16, 6).Value = 0
$ws.Cells.Item(16, 7).Value = 11.68
$ws.Cells.Item(16, 8).Value = "ESE"
$ws.Cells.Item(16, 9).Value = 0
$ws.Cells.Item(16, 10).Value = "20:00:21"
